$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46 (pushes existing rows 46-68 down to 47-69,
# carrying forward formatting from the row above, matching the workbook's
# weekly "insert newest record near the top of this block" pattern).
$ws.Rows(46).Insert()

# Populate the newly inserted row 46 with the new weekly price record.
$ws.Range("A46").Value = 1
$ws.Range("B46").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C46").Value = "Arica y Parinacota"
$ws.Range("D46").Value = 44873
$ws.Range("E46").Value = 15
$ws.Range("F46").Value = 100112031
$ws.Range("G46").Value = "Poroto verde"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 1200
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 900
$ws.Range("M46").Value = 850
$ws.Range("N46").Value = "$/kilo"
$ws.Range("O46").Value = "Región de Arica y Parinacota"
$ws.Range("P46").Value = 850
$ws.Range("Q46").Value = 1
$ws.Range("R46").Value = "Hortaliza"
